$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing text (inline-string) formatting of the data range so
# that writing new numeric-looking values does not silently convert the
# cells to the Number type.
$ws.Range("C2:F37").NumberFormat = "@"

$ws.Range("C2").Value = "100"
$ws.Range("D2").Value = "200"
$ws.Range("E2").Value = "200"
$ws.Range("F2").Value = "50"
$ws.Range("C4").Value = "100"
$ws.Range("D4").Value = "50"
$ws.Range("C5").Value = "10"
$ws.Range("D5").Value = "10"
$ws.Range("E5").Value = "5"
$ws.Range("F5").Value = "10"
$ws.Range("C6").Value = "300"
$ws.Range("F6").Value = "100"
$ws.Range("C8").Value = "0.0001"
$ws.Range("F8").Value = "1e-06"
$ws.Range("C11").Value = "0.1"
$ws.Range("C12").Value = "0.2"
$ws.Range("D12").Value = "0.2"
$ws.Range("E12").Value = "0.2"
$ws.Range("C13").Value = "200"
$ws.Range("D13").Value = "200"
$ws.Range("E13").Value = "200"
$ws.Range("D15").Value = "0.5"
$ws.Range("E15").Value = "0.75"
$ws.Range("C16").Value = "1000"
$ws.Range("D16").Value = "200"
$ws.Range("E16").Value = "100"
$ws.Range("C17").Value = "10"
$ws.Range("D17").Value = "100"
$ws.Range("E17").Value = "100"
$ws.Range("C18").Value = "4"
$ws.Range("C19").Value = "0.01"
$ws.Range("C20").Value = "28"
$ws.Range("D20").Value = "18"
$ws.Range("E20").Value = "28"
$ws.Range("F20").Value = "8"
$ws.Range("C21").Value = "100"
$ws.Range("D21").Value = "100"
$ws.Range("E21").Value = "200"
$ws.Range("E22").Value = "4"
$ws.Range("F22").Value = "5"
$ws.Range("C23").Value = "0.05"
$ws.Range("E23").Value = "0.05"
$ws.Range("E24").Value = "10"
$ws.Range("E25").Value = "0.1"
$ws.Range("D26").Value = "[8]"
$ws.Range("C27").Value = "600"
$ws.Range("D27").Value = "600"
$ws.Range("F28").Value = "10"
$ws.Range("C29").Value = "4"
$ws.Range("D29").Value = "4"
$ws.Range("E29").Value = "5"
$ws.Range("D30").Value = "6"
$ws.Range("C31").Value = "140"
$ws.Range("D31").Value = "120"
$ws.Range("E31").Value = "100"
$ws.Range("F31").Value = "120"
$ws.Range("C32").Value = "0.001"
$ws.Range("E32").Value = "0.001"
$ws.Range("C33").Value = "1000"
$ws.Range("E33").Value = "1000"
$ws.Range("C34").Value = "0.1"
$ws.Range("D34").Value = "0.1"
$ws.Range("F34").Value = "0.1"
$ws.Range("D35").Value = "800"
$ws.Range("E35").Value = "600"
$ws.Range("F35").Value = "600"
$ws.Range("C36").Value = "5"
$ws.Range("F36").Value = "4"
$ws.Range("C37").Value = "0.02"
$ws.Range("E37").Value = "0.01"
$ws.Range("F37").Value = "0.01"

